$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" '43.084.71'
Set-TextValue "E2" '  +1.24%  '

# Row 3
Set-TextValue "D3" '2.269.63'
Set-TextValue "E3" '  +1.46%  '

# Row 4
Set-TextValue "E4" '  +0.14%  '

# Row 5
Set-TextValue "D5" '113.23'
Set-TextValue "E5" '  -1.79%  '

# Row 6
Set-TextValue "D6" '303.41'
Set-TextValue "E6" '  +7.63%  '

# Row 7
Set-TextValue "D7" '0.634'
Set-TextValue "E7" '  +0.97%  '

# Row 8
Set-TextValue "E8" '  -0.16%  '

# Row 9
Set-TextValue "D9" '0.617'
Set-TextValue "E9" '  +0.89%  '

# Row 10
Set-TextValue "D10" '44.50'
Set-TextValue "E10" '  -4.88%  '

# Row 11
Set-TextValue "D11" '0.0930'
Set-TextValue "E11" '  +0.11%  '

# Row 12
Set-TextValue "D12" '54.77'
Set-TextValue "E12" '  +1.11%  '

# Row 13
Set-TextValue "D13" '8.96'
Set-TextValue "E13" '  -2.21%  '

# Row 14
Set-TextValue "E14" '  +20.99%  '

# Row 15
Set-TextValue "E15" '  -0.32%  '

# Row 16
Set-TextValue "D16" '15.46'
Set-TextValue "E16" '  +1.09%  '

# Row 17
Set-TextValue "D17" '2.608.64'
Set-TextValue "E17" '  +1.39%  '

# Row 18
Set-TextValue "D18" '2.316.18'
Set-TextValue "E18" '  +3.73%  '

# Row 19
Set-TextValue "D19" '43.032.18'
Set-TextValue "E19" '  +0.69%  '

# Row 20
Set-TextValue "D20" '0.0000108'
Set-TextValue "E20" '  +0.28%  '

# Row 21
Set-TextValue "D21" '7.25'
Set-TextValue "E21" '  +5.68%  '

# Row 22
Set-TextValue "D22" '75.39'
Set-TextValue "E22" '  +4.42%  '

# Row 23
Set-TextValue "D23" '3.58'
Set-TextValue "E23" '  +15.67%  '

# Row 24
Set-TextValue "D24" '258.64'
Set-TextValue "E24" '  +11.55%  '

# Row 25
Set-TextValue "D25" '2.44'
Set-TextValue "E25" '  +4.35%  '

# Row 26
Set-TextValue "D26" '9.02'
Set-TextValue "E26" '  -3.29%  '

# Row 27
Set-TextValue "B27" 'Cosmos'
Set-TextValue "C27" 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D27" '11.67'
Set-TextValue "E27" '  -3.26%  '

# Row 28
Set-TextValue "B28" 'Dai'
Set-TextValue "C28" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue "D28" '1.00'
Set-TextValue "E28" '  -0.19%  '

# Row 29
Set-TextValue "E29" '  -0.39%  '

# Row 30
Set-TextValue "D30" '38.27'
Set-TextValue "E30" '  -5.22%  '

# Row 31
Set-TextValue "D31" '22.41'
Set-TextValue "E31" '  +5.94%  '

# Row 32
Set-TextValue "D32" '175.17'
Set-TextValue "E32" '  +0.74%  '

# Row 33
Set-TextValue "D33" '3.19'
Set-TextValue "E33" '  -2.83%  '

# Row 34
Set-TextValue "D34" '0.0896'
Set-TextValue "E34" '  -0.29%  '

# Row 35
Set-TextValue "D35" '5.72'
Set-TextValue "E35" '  +2.45%  '

# Row 36
Set-TextValue "D36" '5.09'
Set-TextValue "E36" '  +9.49%  '

# Row 37
Set-TextValue "B37" 'NEARProtocol'
Set-TextValue "C37" 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D37" '4.27'
Set-TextValue "E37" '  -4.92%  '

# Row 38
Set-TextValue "B38" 'Stellar'
Set-TextValue "C38" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D38" '0.129'
Set-TextValue "E38" '  +0.38%  '

# Row 39
Set-TextValue "D39" '0.0378'
Set-TextValue "E39" '  +1.76%  '

# Row 40
Set-TextValue "D40" '0.105'
Set-TextValue "E40" '  -1.33%  '

# Row 41
Set-TextValue "D41" '2.46'
Set-TextValue "E41" '  -5.58%  '

# Row 42
Set-TextValue "D42" '72.35'
Set-TextValue "E42" '  +0.36%  '

# Row 43
Set-TextValue "D43" '0.233'
Set-TextValue "E43" '  -0.68%  '

# Row 44
Set-TextValue "E44" '  -0.22%  '

# Row 45
Set-TextValue "D45" '12.63'
Set-TextValue "E45" '  -6.69%  '

# Row 46
Set-TextValue "D46" '1.35'
Set-TextValue "E46" '  +0.65%  '

# Row 47
Set-TextValue "D47" '5.56'
Set-TextValue "E47" '  +0.23%  '

# Row 48
Set-TextValue "D48" '107.93'
Set-TextValue "E48" '  +7.17%  '

# Row 49
Set-TextValue "D49" '1.31'
Set-TextValue "E49" '  +1.77%  '

# Row 50
Set-TextValue "E50" '  +2.69%  '

# Row 51
Set-TextValue "D51" '72.99'
Set-TextValue "E51" '  +3.34%  '
